$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (585) down to the new rows (586-601)
$ws.Range("A585:V585").Copy()
$ws.Range("A586:V601").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 586
$ws.Range("A586").Value = 'Entrainement'
$ws.Range("B586").Value = 45931
$ws.Range("C586").Value = 'Global'
$ws.Range("D586").Value = 'J-3'
$ws.Range("E586").Value = 'Naim Dhib'
$ws.Range("F586").Value = 'center midfield'
$ws.Range("G586").Value = '01:38:54'
$ws.Range("H586").Value = 6.72
$ws.Range("I586").Value = 0.94
$ws.Range("J586").Value = 5.77
$ws.Range("K586").Value = 0.55
$ws.Range("L586").Value = 0.32
$ws.Range("M586").Value = 0.08
$ws.Range("N586").Value = 0
$ws.Range("O586").Value = 11
$ws.Range("P586").Value = 3.49
$ws.Range("Q586").Value = 27.9
$ws.Range("R586").Value = 5.27
$ws.Range("S586").Value = 32
$ws.Range("T586").Value = 3
$ws.Range("U586").Value = 25
$ws.Range("V586").Value = 8

# Row 587
$ws.Range("A587").Value = 'Entrainement'
$ws.Range("B587").Value = 45931
$ws.Range("C587").Value = 'Global'
$ws.Range("D587").Value = 'J-3'
$ws.Range("E587").Value = 'Amir Etien'
$ws.Range("F587").Value = 'right forward'
$ws.Range("G587").Value = '01:36:13'
$ws.Range("H587").Value = 6.22
$ws.Range("I587").Value = 0.99
$ws.Range("J587").Value = 5.22
$ws.Range("K587").Value = 0.51
$ws.Range("L587").Value = 0.31
$ws.Range("M587").Value = 0.16
$ws.Range("N587").Value = 0.03
$ws.Range("O587").Value = 13
$ws.Range("P587").Value = 3.4
$ws.Range("Q587").Value = 35.4
$ws.Range("R587").Value = 5.77
$ws.Range("S587").Value = 49
$ws.Range("T587").Value = 15
$ws.Range("U587").Value = 31
$ws.Range("V587").Value = 13

# Row 588
$ws.Range("A588").Value = 'Entrainement'
$ws.Range("B588").Value = 45931
$ws.Range("C588").Value = 'Global'
$ws.Range("D588").Value = 'J-3'
$ws.Range("E588").Value = 'Karahali Souaré'
$ws.Range("F588").Value = 'right forward'
$ws.Range("G588").Value = '01:38:54'
$ws.Range("H588").Value = 6.62
$ws.Range("I588").Value = 1.19
$ws.Range("J588").Value = 5.41
$ws.Range("K588").Value = 0.56
$ws.Range("L588").Value = 0.35
$ws.Range("M588").Value = 0.25
$ws.Range("N588").Value = 0.05
$ws.Range("O588").Value = 21
$ws.Range("P588").Value = 3.65
$ws.Range("Q588").Value = 33.2
$ws.Range("R588").Value = 5.62
$ws.Range("S588").Value = 59
$ws.Range("T588").Value = 15
$ws.Range("U588").Value = 44
$ws.Range("V588").Value = 24

# Row 589
$ws.Range("A589").Value = 'Entrainement'
$ws.Range("B589").Value = 45931
$ws.Range("C589").Value = 'Global'
$ws.Range("D589").Value = 'J-3'
$ws.Range("E589").Value = 'Amir Kherrab'
$ws.Range("F589").Value = 'center midfield'
$ws.Range("G589").Value = '00:45:48'
$ws.Range("H589").Value = 2.86
$ws.Range("I589").Value = 0.28
$ws.Range("J589").Value = 2.57
$ws.Range("K589").Value = 0.19
$ws.Range("L589").Value = 0.08
$ws.Range("M589").Value = 0.01
$ws.Range("N589").Value = 0
$ws.Range("O589").Value = 1
$ws.Range("P589").Value = 3.69
$ws.Range("Q589").Value = 25.77
$ws.Range("R589").Value = 3.85
$ws.Range("S589").Value = 12
$ws.Range("T589").Value = 0
$ws.Range("U589").Value = 7
$ws.Range("V589").Value = 2

# Row 590
$ws.Range("A590").Value = 'Entrainement'
$ws.Range("B590").Value = 45931
$ws.Range("C590").Value = 'Global'
$ws.Range("D590").Value = 'J-3'
$ws.Range("E590").Value = 'Kamal Bafounta'
$ws.Range("F590").Value = 'center midfield'
$ws.Range("G590").Value = '01:40:35'
$ws.Range("H590").Value = 7.03
$ws.Range("I590").Value = 0.83
$ws.Range("J590").Value = 6.18
$ws.Range("K590").Value = 0.57
$ws.Range("L590").Value = 0.2
$ws.Range("M590").Value = 0.08
$ws.Range("N590").Value = 0
$ws.Range("O590").Value = 3
$ws.Range("P590").Value = 4.11
$ws.Range("Q590").Value = 29.96
$ws.Range("R590").Value = 4.43
$ws.Range("S590").Value = 42
$ws.Range("T590").Value = 6
$ws.Range("U590").Value = 27
$ws.Range("V590").Value = 6

# Row 591
$ws.Range("A591").Value = 'Entrainement'
$ws.Range("B591").Value = 45931
$ws.Range("C591").Value = 'Global'
$ws.Range("D591").Value = 'J-3'
$ws.Range("E591").Value = 'Yoann Martelat'
$ws.Range("F591").Value = 'center midfield'
$ws.Range("G591").Value = '01:38:39'
$ws.Range("H591").Value = 7.33
$ws.Range("I591").Value = 1.04
$ws.Range("J591").Value = 6.28
$ws.Range("K591").Value = 0.68
$ws.Range("L591").Value = 0.32
$ws.Range("M591").Value = 0.05
$ws.Range("N591").Value = 0.01
$ws.Range("O591").Value = 3
$ws.Range("P591").Value = 4.35
$ws.Range("Q591").Value = 30.34
$ws.Range("R591").Value = 3.98
$ws.Range("S591").Value = 17
$ws.Range("T591").Value = 0
$ws.Range("U591").Value = 11
$ws.Range("V591").Value = 2

# Row 592
$ws.Range("A592").Value = 'Entrainement'
$ws.Range("B592").Value = 45931
$ws.Range("C592").Value = 'Global'
$ws.Range("D592").Value = 'J-3'
$ws.Range("E592").Value = 'Malik Boussaid'
$ws.Range("F592").Value = 'right back'
$ws.Range("G592").Value = '01:38:39'
$ws.Range("H592").Value = 7.71
$ws.Range("I592").Value = 1.37
$ws.Range("J592").Value = 6.32
$ws.Range("K592").Value = 0.69
$ws.Range("L592").Value = 0.37
$ws.Range("M592").Value = 0.28
$ws.Range("N592").Value = 0.05
$ws.Range("O592").Value = 12
$ws.Range("P592").Value = 4.05
$ws.Range("Q592").Value = 31.57
$ws.Range("R592").Value = 4.76
$ws.Range("S592").Value = 35
$ws.Range("T592").Value = 6
$ws.Range("U592").Value = 31
$ws.Range("V592").Value = 13

# Row 593
$ws.Range("A593").Value = 'Entrainement'
$ws.Range("B593").Value = 45931
$ws.Range("C593").Value = 'Global'
$ws.Range("D593").Value = 'J-3'
$ws.Range("E593").Value = 'Sofiane Belle'
$ws.Range("F593").Value = 'left forward'
$ws.Range("G593").Value = '01:35:57'
$ws.Range("H593").Value = 5.72
$ws.Range("I593").Value = 0.83
$ws.Range("J593").Value = 4.88
$ws.Range("K593").Value = 0.41
$ws.Range("L593").Value = 0.28
$ws.Range("M593").Value = 0.15
$ws.Range("N593").Value = 0
$ws.Range("O593").Value = 11
$ws.Range("P593").Value = 3.37
$ws.Range("Q593").Value = 29.44
$ws.Range("R593").Value = 4.27
$ws.Range("S593").Value = 9
$ws.Range("T593").Value = 1
$ws.Range("U593").Value = 19
$ws.Range("V593").Value = 5

# Row 594
$ws.Range("A594").Value = 'Entrainement'
$ws.Range("B594").Value = 45931
$ws.Range("C594").Value = 'Global'
$ws.Range("D594").Value = 'J-3'
$ws.Range("E594").Value = 'Jeremie Laurent'
$ws.Range("F594").Value = 'left forward'
$ws.Range("G594").Value = '01:37:22'
$ws.Range("H594").Value = 6.88
$ws.Range("I594").Value = 1.48
$ws.Range("J594").Value = 5.39
$ws.Range("K594").Value = 0.77
$ws.Range("L594").Value = 0.44
$ws.Range("M594").Value = 0.26
$ws.Range("N594").Value = 0.02
$ws.Range("O594").Value = 16
$ws.Range("P594").Value = 4.14
$ws.Range("Q594").Value = 31.54
$ws.Range("R594").Value = 4.71
$ws.Range("S594").Value = 34
$ws.Range("T594").Value = 13
$ws.Range("U594").Value = 20
$ws.Range("V594").Value = 14

# Row 595
$ws.Range("A595").Value = 'Entrainement'
$ws.Range("B595").Value = 45931
$ws.Range("C595").Value = 'Global'
$ws.Range("D595").Value = 'J-3'
$ws.Range("E595").Value = 'Mattheo Haon'
$ws.Range("F595").Value = 'right back'
$ws.Range("G595").Value = '01:39:40'
$ws.Range("H595").Value = 7.84
$ws.Range("I595").Value = 1.23
$ws.Range("J595").Value = 6.6
$ws.Range("K595").Value = 0.76
$ws.Range("L595").Value = 0.27
$ws.Range("M595").Value = 0.17
$ws.Range("N595").Value = 0.05
$ws.Range("O595").Value = 10
$ws.Range("P595").Value = 4.66
$ws.Range("Q595").Value = 31.18
$ws.Range("R595").Value = 4.94
$ws.Range("S595").Value = 28
$ws.Range("T595").Value = 7
$ws.Range("U595").Value = 25
$ws.Range("V595").Value = 7

# Row 596
$ws.Range("A596").Value = 'Entrainement'
$ws.Range("B596").Value = 45931
$ws.Range("C596").Value = 'Global'
$ws.Range("D596").Value = 'J-3'
$ws.Range("E596").Value = 'Ilan Ihaddadene'
$ws.Range("F596").Value = 'center midfield'
$ws.Range("G596").Value = '01:38:46'
$ws.Range("H596").Value = 7.67
$ws.Range("I596").Value = 1.22
$ws.Range("J596").Value = 6.44
$ws.Range("K596").Value = 0.91
$ws.Range("L596").Value = 0.29
$ws.Range("M596").Value = 0.03
$ws.Range("N596").Value = 0
$ws.Range("O596").Value = 3
$ws.Range("P596").Value = 4.56
$ws.Range("Q596").Value = 26.94
$ws.Range("R596").Value = 5.05
$ws.Range("S596").Value = 37
$ws.Range("T596").Value = 5
$ws.Range("U596").Value = 21
$ws.Range("V596").Value = 4

# Row 597
$ws.Range("A597").Value = 'Entrainement'
$ws.Range("B597").Value = 45931
$ws.Range("C597").Value = 'Global'
$ws.Range("D597").Value = 'J-3'
$ws.Range("E597").Value = 'Karim Belmahi'
$ws.Range("F597").Value = 'left forward'
$ws.Range("G597").Value = '01:37:45'
$ws.Range("H597").Value = 7.47
$ws.Range("I597").Value = 1.19
$ws.Range("J597").Value = 6.25
$ws.Range("K597").Value = 0.73
$ws.Range("L597").Value = 0.31
$ws.Range("M597").Value = 0.14
$ws.Range("N597").Value = 0.05
$ws.Range("O597").Value = 14
$ws.Range("P597").Value = 4.32
$ws.Range("Q597").Value = 33.23
$ws.Range("R597").Value = 5.75
$ws.Range("S597").Value = 58
$ws.Range("T597").Value = 20
$ws.Range("U597").Value = 47
$ws.Range("V597").Value = 24

# Row 598
$ws.Range("A598").Value = 'Entrainement'
$ws.Range("B598").Value = 45931
$ws.Range("C598").Value = 'Global'
$ws.Range("D598").Value = 'J-3'
$ws.Range("E598").Value = 'Omar Benyounes'
$ws.Range("F598").Value = 'center midfield'
$ws.Range("G598").Value = '01:39:02'
$ws.Range("H598").Value = 7.65
$ws.Range("I598").Value = 1.04
$ws.Range("J598").Value = 6.59
$ws.Range("K598").Value = 0.68
$ws.Range("L598").Value = 0.26
$ws.Range("M598").Value = 0.11
$ws.Range("N598").Value = 0
$ws.Range("O598").Value = 7
$ws.Range("P598").Value = 4.54
$ws.Range("Q598").Value = 27.38
$ws.Range("R598").Value = 4.5
$ws.Range("S598").Value = 35
$ws.Range("T598").Value = 4
$ws.Range("U598").Value = 27
$ws.Range("V598").Value = 5

# Row 599
$ws.Range("A599").Value = 'Entrainement'
$ws.Range("B599").Value = 45931
$ws.Range("C599").Value = 'Global'
$ws.Range("D599").Value = 'J-3'
$ws.Range("E599").Value = 'Hedi Nasri'
$ws.Range("F599").Value = 'right back'
$ws.Range("G599").Value = '01:40:20'
$ws.Range("H599").Value = 6.15
$ws.Range("I599").Value = 0.86
$ws.Range("J599").Value = 5.28
$ws.Range("K599").Value = 0.58
$ws.Range("L599").Value = 0.24
$ws.Range("M599").Value = 0.06
$ws.Range("N599").Value = 0
$ws.Range("O599").Value = 9
$ws.Range("P599").Value = 3.48
$ws.Range("Q599").Value = 28.02
$ws.Range("R599").Value = 5.7
$ws.Range("S599").Value = 47
$ws.Range("T599").Value = 7
$ws.Range("U599").Value = 27
$ws.Range("V599").Value = 16

# Row 600
$ws.Range("A600").Value = 'Entrainement'
$ws.Range("B600").Value = 45931
$ws.Range("C600").Value = 'Global'
$ws.Range("D600").Value = 'J-3'
$ws.Range("E600").Value = 'Levy Ndoutoume'
$ws.Range("F600").Value = 'left back'
$ws.Range("G600").Value = '01:39:17'
$ws.Range("H600").Value = 6.88
$ws.Range("I600").Value = 0.84
$ws.Range("J600").Value = 6.02
$ws.Range("K600").Value = 0.48
$ws.Range("L600").Value = 0.3
$ws.Range("M600").Value = 0.07
$ws.Range("N600").Value = 0
$ws.Range("O600").Value = 7
$ws.Range("P600").Value = 3.77
$ws.Range("Q600").Value = 29.99
$ws.Range("R600").Value = 4.5
$ws.Range("S600").Value = 42
$ws.Range("T600").Value = 4
$ws.Range("U600").Value = 33
$ws.Range("V600").Value = 9

# Row 601
$ws.Range("A601").Value = 'Entrainement'
$ws.Range("B601").Value = 45931
$ws.Range("C601").Value = 'Global'
$ws.Range("D601").Value = 'J-3'
$ws.Range("E601").Value = 'Emmanuel Valey'
$ws.Range("F601").Value = 'left forward'
$ws.Range("G601").Value = '01:39:09'
$ws.Range("H601").Value = 7.46
$ws.Range("I601").Value = 1.23
$ws.Range("J601").Value = 6.21
$ws.Range("K601").Value = 0.57
$ws.Range("L601").Value = 0.42
$ws.Range("M601").Value = 0.2
$ws.Range("N601").Value = 0.06
$ws.Range("O601").Value = 19
$ws.Range("P601").Value = 4.04
$ws.Range("Q601").Value = 32.27
$ws.Range("R601").Value = 5.06
$ws.Range("S601").Value = 43
$ws.Range("T601").Value = 14
$ws.Range("U601").Value = 26
$ws.Range("V601").Value = 17
# Update the view: selection moves to B605 (topLeftCell A567 is not representable via this COM surface)
$excel.Goto($ws.Range("B605"), $false)
$ws.Range("B605").Select()
